$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value2 = 0.7376376588883126
$ws.Range("J2").Value2 = 0.7376376588883125
$ws.Range("M2").Value2 = 9.506955
$ws.Range("N2").Value2 = 28.520865
$ws.Range("O2").Value2 = 0.2691767467564006
$ws.Range("P2").Value2 = 0.2691767467564006
$ws.Range("Q2").Value2 = 5.53851114014
$ws.Range("R2").Value2 = 49.84660026126
$ws.Range("S2").Value2 = 0.1985549053045635
$ws.Range("T2").Value2 = 0.1985549053045635

$ws.Range("I3").Value2 = 0.7376376588883126
$ws.Range("J3").Value2 = 0.7376376588883125
$ws.Range("O3").Value2 = 0.5092171458273274
$ws.Range("P3").Value2 = 0.5092171458273272
$ws.Range("S3").Value2 = 0.3756177433138582
$ws.Range("T3").Value2 = 0.3756177433138581

$ws.Range("I4").Value2 = 0.7376376588883126
$ws.Range("J4").Value2 = 0.7376376588883125
$ws.Range("M4").Value2 = 7.826824999999999
$ws.Range("N4").Value2 = 23.480475
$ws.Range("O4").Value2 = 0.2216061074162721
$ws.Range("P4").Value2 = 0.2216061074162721
$ws.Range("Q4").Value2 = 4.559709965433333
$ws.Range("R4").Value2 = 41.0373896889
$ws.Range("S4").Value2 = 0.1634650102698909
$ws.Range("T4").Value2 = 0.1634650102698909

$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 0.6666666666666666
$ws.Range("G5").Value2 = 0.2072096666666667
$ws.Range("H5").Value2 = 0.621629
$ws.Range("I5").Value2 = 0.2623623411116874
$ws.Range("J5").Value2 = 0.2623623411116874
$ws.Range("M5").Value2 = 9.506955
$ws.Range("N5").Value2 = 28.520865
$ws.Range("O5").Value2 = 0.2691767467564006
$ws.Range("P5").Value2 = 0.2691767467564006
$ws.Range("Q5").Value2 = 1.969932976565
$ws.Range("R5").Value2 = 17.729396789085
$ws.Range("S5").Value2 = 0.07062184145183709
$ws.Range("T5").Value2 = 0.07062184145183709

$ws.Range("E6").Value2 = 2
$ws.Range("F6").Value2 = 0.6666666666666666
$ws.Range("G6").Value2 = 0.2072096666666667
$ws.Range("H6").Value2 = 0.621629
$ws.Range("I6").Value2 = 0.2623623411116874
$ws.Range("J6").Value2 = 0.2623623411116874
$ws.Range("O6").Value2 = 0.5092171458273274
$ws.Range("P6").Value2 = 0.5092171458273272
$ws.Range("Q6").Value2 = 3.726635602388666
$ws.Range("R6").Value2 = 33.539720421498
$ws.Range("S6").Value2 = 0.1335994025134691
$ws.Range("T6").Value2 = 0.1335994025134691

$ws.Range("E7").Value2 = 2
$ws.Range("F7").Value2 = 0.6666666666666666
$ws.Range("G7").Value2 = 0.2072096666666667
$ws.Range("H7").Value2 = 0.621629
$ws.Range("I7").Value2 = 0.2623623411116874
$ws.Range("J7").Value2 = 0.2623623411116874
$ws.Range("M7").Value2 = 7.826824999999999
$ws.Range("N7").Value2 = 23.480475
$ws.Range("O7").Value2 = 0.2216061074162721
$ws.Range("P7").Value2 = 0.2216061074162721
$ws.Range("Q7").Value2 = 1.621793799308333
$ws.Range("R7").Value2 = 14.596144193775
$ws.Range("S7").Value2 = 0.05814109714638124
$ws.Range("T7").Value2 = 0.05814109714638123

